$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New summary rows added below the data table ---

# Row 12: average of the k column (J), bold
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"

# Row 14-17: labelled summary statistics in column A/B
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

# --- Formatting ---

# J12: bold
$ws.Range("J12").Font.Bold = $true

# B14:B17: bold, 12pt, vertically centered
$summaryRange = $ws.Range("B14:B17")
$summaryRange.Font.Bold = $true
$summaryRange.Font.Size = 12
$summaryRange.VerticalAlignment = -4108

# --- Page setup (portrait, paper size 9 = A4) ---
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# --- Selection matches the author's last saved cursor position ---
$null = $ws.Range("J12").Select()
